# Fix excel conversion script: annotation/plot files were moved into a new
# "conservation_analysis" subfolder, and a "multi_level_plot" png column
# (M) was added for the pairk_aln_needleman... method. The json-file
# columns (L, Z) are no longer turned into hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a copy of the "link" cell formatting (font/style) used by
# G4/U4/AA4/.. in an unused scratch cell (well outside the used range)
# so it can be re-applied later, since rebuilding the Hyperlinks
# collection resets cell formatting.
$ws.Range("G4").Copy()
$ws.Range("AI1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 1. Update displayed paths for cells whose files moved into the new
#    conservation_analysis/annotations/ subdirectory.
# ---------------------------------------------------------------------
$ws.Range("G4").Value  = "conservation_analysis/annotations/2-9606_0002f40-Vertebrata_aln_slice.html"
$ws.Range("U4").Value  = "conservation_analysis/annotations/2-9606_0002f40-Vertebrata_aln_slice.html"
$ws.Range("AA4").Value = "conservation_analysis/annotations/2-9606_0002f40-aln_property_entropy_multilevel_plot.png"

$ws.Range("G5").Value  = "conservation_analysis/annotations/3-9606_0002f40-Vertebrata_aln_slice.html"
$ws.Range("U5").Value  = "conservation_analysis/annotations/3-9606_0002f40-Vertebrata_aln_slice.html"
$ws.Range("AA5").Value = "conservation_analysis/annotations/3-9606_0002f40-aln_property_entropy_multilevel_plot.png"

# ---------------------------------------------------------------------
# 2. The json-file columns (L, Z) are no longer rendered as hyperlinks
#    in the fixed script, so drop their hyperlink styling.
# ---------------------------------------------------------------------
$ws.Range("L4").ClearFormats()
$ws.Range("Z4").ClearFormats()
$ws.Range("L5").ClearFormats()
$ws.Range("Z5").ClearFormats()

# ---------------------------------------------------------------------
# 3. Add the new "multi_level_plot" png column (M) for the
#    pairk_aln_needleman... method, styled like the other link cells.
# ---------------------------------------------------------------------
$ws.Range("AI1").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = "conservation_analysis/annotations/2-9606_0002f40-pairk_aln_needleman_lf5_rf5_edssmat50_multilevel_plot.png"

$ws.Range("AI1").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = "conservation_analysis/annotations/3-9606_0002f40-pairk_aln_needleman_lf5_rf5_edssmat50_multilevel_plot.png"

# ---------------------------------------------------------------------
# 4. Rebuild the hyperlinks: drop the ones for L4/Z4/L5/Z5 (json files),
#    keep the others pointing at their (moved) files, and add new ones
#    for the new M4/M5 plot cells.
# ---------------------------------------------------------------------
$base = "file:///Users/jackson/Dropbox (MIT)/work/07-SLiM_bioinformatics/05-conservation_pipeline/examples/table_annotation/"

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G4"),  $base + "conservation_analysis/annotations/2-9606_0002f40-Vertebrata_aln_slice.html")
$ws.Hyperlinks.Add($ws.Range("M4"),  $base + "conservation_analysis/annotations/2-9606_0002f40-pairk_aln_needleman_lf5_rf5_edssmat50_multilevel_plot.png")
$ws.Hyperlinks.Add($ws.Range("U4"),  $base + "conservation_analysis/annotations/2-9606_0002f40-Vertebrata_aln_slice.html")
$ws.Hyperlinks.Add($ws.Range("AA4"), $base + "conservation_analysis/annotations/2-9606_0002f40-aln_property_entropy_multilevel_plot.png")
$ws.Hyperlinks.Add($ws.Range("G5"),  $base + "conservation_analysis/annotations/3-9606_0002f40-Vertebrata_aln_slice.html")
$ws.Hyperlinks.Add($ws.Range("M5"),  $base + "conservation_analysis/annotations/3-9606_0002f40-pairk_aln_needleman_lf5_rf5_edssmat50_multilevel_plot.png")
$ws.Hyperlinks.Add($ws.Range("U5"),  $base + "conservation_analysis/annotations/3-9606_0002f40-Vertebrata_aln_slice.html")
$ws.Hyperlinks.Add($ws.Range("AA5"), $base + "conservation_analysis/annotations/3-9606_0002f40-aln_property_entropy_multilevel_plot.png")

# Re-adding hyperlinks resets the cell formatting to the built-in
# "Hyperlink" style, so restore the original link-cell formatting that
# was stashed above.
$ws.Range("AI1").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("U4").PasteSpecial(-4122)
$ws.Range("AA4").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("U5").PasteSpecial(-4122)
$ws.Range("AA5").PasteSpecial(-4122)

# Remove the scratch column entirely so it leaves no trace in the sheet.
$ws.Range("AI1").EntireColumn.Delete()

Write-Output "done"
